# Update cryptocurrency price/volume data (and restore correct Polkadot/Chainlink
# and NEARProtocol/RenderToken row ordering) per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.802.01"
$ws.Range("E2").Value = "  +2.97%  "

$ws.Range("D3").Value = "1.880.64"
$ws.Range("E3").Value = "  +3.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.54%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.45%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4677"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.03%  "

$ws.Range("E8").Value = "  +2.37%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07926"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.81%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9800"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.40%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.33"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.17%  "

$ws.Range("D12").Value = "1.949.53"
$ws.Range("E12").Value = "  +6.85%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.742"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.84%  "

$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.017"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06956"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.90%  "

$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("E18").Value = "  +1.63%  "

$ws.Range("E19").Value = "  +2.17%  "

$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("D21").Value = "28.810.54"
$ws.Range("E21").Value = "  +2.91%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.344"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.10"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.27%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.118"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.02%  "

$ws.Range("D25").Value = "2.079.16"
$ws.Range("E25").Value = "  +1.99%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.83"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.749"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.001"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.71%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9396"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.320"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.85%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.357"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.36%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.05918"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.05%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02122"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.158"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.40%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.889"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5728"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.88%  "

$ws.Range("E41").Value = "  +2.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.00"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.07327"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.76%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.90"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5340"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.98%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.155"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.26%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.116"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.51%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.848"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.58%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "114.08"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.374"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.90%  "

$ws.Range("E51").Value = "  +0.65%  "
